$d = $word.ActiveDocument
$full = $d.Content.WordOpenXML
$idx = $full.IndexOf("<w:document ")
$idx2 = $full.IndexOf(">", $idx)
Write-Output ($idx2 - $idx)
Write-Output $full.Substring($idx, $idx2 - $idx + 1)
